$p = $ppt.ActivePresentation

# Slide 1: "Machine Learning Fundamentals" - bullet content expanded to 3 lines
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "What ML is`rWhy it matters`rWhere it is used"

# Slide 2: "What is Machine Learning?" - bullet content expanded to 3 lines
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "Algorithms that learn patterns`rImprove with data`rMake predictions or decisions"

# Slide 3: "Types of Machine Learning" - bullet content expanded to 3 lines
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Supervised learning`rUnsupervised learning`rReinforcement learning"

# Slide 4: "Supervised Learning" - bullet content expanded to 3 lines
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Uses labelled data`rRegression & classification`rExamples: spam detection, price prediction"

# Slide 5: "Unsupervised Learning" - bullet content expanded to 3 lines
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "No labelled outputs`rPattern discovery`rExamples: clustering customers"

# Slide 6: "Reinforcement Learning" - bullet content expanded to 3 lines
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Agent learns by rewards`rTrial and error`rExamples: games, robotics"

# Slide 7: "Common Algorithms" - bullet content expanded to 4 lines
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(2).TextFrame.TextRange.Text = "Linear regression`rLogistic regression`rDecision trees`rK-means"

# Slide 8: Title renamed, bullet content expanded to 4 lines
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Model Training Process"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "Collect data`rTrain model`rEvaluate performance`rDeploy"

# Slide 9: Title renamed, bullet content replaced with 3 lines
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Overfitting vs Underfitting"
$s9.Shapes.Item(2).TextFrame.TextRange.Text = "Overfitting: too complex`rUnderfitting: too simple`rBias" + [char]0x2013 + "variance tradeoff"

# Slide 10: Title renamed, bullet content expanded to 3 lines
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Key Takeaways"
$s10.Shapes.Item(2).TextFrame.TextRange.Text = "ML learns from data`rDifferent learning types`rUsed across industries"
